# "New UI changes in Candidate Lobby"
#
# 1) AMSIN sheet: append three new interview-history rows (70-72) after the
#    existing data in A1:G69, extending the used range to A1:G72.
# 2) AMS sheet: normalize the formatting on row 48 (the most recently added
#    row) so it matches the style used by the rest of the data rows, and
#    correct the run-time serial value in B48.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. AMSIN - add rows 70, 71, 72
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$newRows = @(
    @{ Row = 70; A = "2023-05-09"; B = 45055.69569077546; C = "176scndhtfx"; D = 155; E = 147; F = 8; G = 4.54 },
    @{ Row = 71; A = "2023-05-11"; B = 45057.75206074074; C = "176fxhh";     D = 155; E = 149; F = 6; G = 4.75 },
    @{ Row = 72; A = "2023-05-12"; B = 45058.76826271369; C = "177fstcycle"; D = 155; E = 153; F = 2; G = 5.03 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Text columns (A = run date, C = sprint name): force literal text so
    # Excel doesn't auto-convert "2023-05-09"-style strings into real dates.
    $wsAmsin.Cells.Item($row, 1).NumberFormat = "@"
    $wsAmsin.Cells.Item($row, 1).Value = $r.A
    $wsAmsin.Cells.Item($row, 3).NumberFormat = "@"
    $wsAmsin.Cells.Item($row, 3).Value = $r.C

    $wsAmsin.Cells.Item($row, 2).Value = $r.B
    $wsAmsin.Cells.Item($row, 4).Value = $r.D
    $wsAmsin.Cells.Item($row, 5).Value = $r.E
    $wsAmsin.Cells.Item($row, 6).Value = $r.F
    $wsAmsin.Cells.Item($row, 7).Value = $r.G

    # Re-apply the same look-and-feel as the rest of the table (the row
    # directly above is already formatted the way every other data row is).
    $wsAmsin.Range("A69:G69").Copy()
    $wsAmsin.Range("A$row`:G$row").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 2. AMS - fix up row 48 formatting + B48 run-time precision
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Cells.Item(48, 2).Value = 45054.54037847222

# Row 47 carries the formatting every other data row uses; bring row 48 in
# line with it (it was pasted in without the usual style).
$wsAms.Range("A47:G47").Copy()
$wsAms.Range("A48:G48").PasteSpecial(-4122)

# PasteSpecial(xlPasteFormats) only carries formatting - restore the actual
# values/content for row 48 since the paste could not have touched them.
$wsAms.Cells.Item(48, 1).NumberFormat = "@"
$wsAms.Cells.Item(48, 1).Value = "2023-05-08"
$wsAms.Cells.Item(48, 2).Value = 45054.54037847222
$wsAms.Cells.Item(48, 3).NumberFormat = "@"
$wsAms.Cells.Item(48, 3).Value = "176htfxtrl"
$wsAms.Cells.Item(48, 4).Value = 155
$wsAms.Cells.Item(48, 5).Value = 155
$wsAms.Cells.Item(48, 6).Value = 0
$wsAms.Cells.Item(48, 7).Value = 3.02
